$d = $word.ActiveDocument
$d.Content.Find.Execute("Clarify the divison of responsibility.", $true, $false, $false, $false, $false, $true, 1, $false, "Clarify the division of responsibility.", 2)
